# Fixed typo in presentation for IETF 116
#
# Slide 3 ("Issues and To do"), shape "Content Placeholder 2":
#   - Paragraph 3 ("Location of tp-to-interface-ref node in the YANG
#     tree...") loses the stray note "Add this information to the
#     presentation for IETF 116." that had leaked into the bullet text,
#     leaving the closing ". " as its own run.
#   - Paragraph 4 ("Which working group will be the best home...") gets
#     its leading "Which " word split into its own run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# --- Paragraph 3: "...Location of tp-to-interface-ref node in the YANG
#     tree. The authors agree, but the change has not yet been
#     implemented. Add this information to the presentation for IETF 116."
$para3 = $tr.Paragraphs(3, 1)

# Drop the stray note sentence, keeping the final period + space.
$note = $para3.Find(". Add this information to the presentation for IETF 116.")
$note.Text = ""

# Re-add the ". " as its own trailing run (matches the author's split).
[void]$para3.InsertAfter(". ")

# --- Paragraph 4: "Which working group will be the best home to progress
#     if-ref-topo-yang and bwa-topo-yang"
$para4 = $tr.Paragraphs(4, 1)

# Split the leading "Which " off into its own run.
$lead = $para4.Find("Which ")
$lead.Text = ""
[void]$para4.InsertBefore("Which ")
